$d = $word.ActiveDocument

$replacements = @(
    @("342÷6=57, 0", "375÷7=53, 4"),
    @("759÷8=94, 7", "786÷8=98, 2"),
    @("632÷5=126, 2", "835÷3=278, 1"),
    @("937÷3=312, 1", "572÷3=190, 2"),
    @("366÷4=91, 2", "105÷8=13, 1"),
    @("223÷9=24, 7", "325÷6=54, 1"),
    @("474÷7=67, 5", "591÷9=65, 6"),
    @("663÷5=132, 3", "756÷4=189, 0"),
    @("266÷3=88, 2", "608÷5=121, 3"),
    @("863÷4=215, 3", "900÷7=128, 4"),
    @("773÷6=128, 5", "224÷9=24, 8"),
    @("209÷4=52, 1", "213÷7=30, 3"),
    @("485÷6=80, 5", "734÷9=81, 5"),
    @("485÷4=121, 1", "918÷7=131, 1"),
    @("408÷6=68, 0", "579÷9=64, 3"),
    @("169÷5=33, 4", "661÷3=220, 1"),
    @("392÷5=78, 2", "253÷6=42, 1"),
    @("935÷8=116, 7", "925÷6=154, 1"),
    @("481÷4=120, 1", "543÷6=90, 3"),
    @("240÷8=30, 0", "305÷3=101, 2"),
    @("146÷4=36, 2", "773÷3=257, 2"),
    @("894÷6=149, 0", "452÷8=56, 4"),
    @("109÷6=18, 1", "592÷3=197, 1"),
    @("228÷8=28, 4", "989÷3=329, 2"),
    @("568÷4=142, 0", "670÷7=95, 5")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
